$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix malformed email addresses (digit suffix moved before the @ sign)
$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"

# Remove the trailing (duplicate/stray) email cells in rows 7 and 8
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Restore the active selection to G4
$ws.Range("G4").Select()
